# Rename the "tagger" model and certain "judged" model names in the
# header row (row 1) of Sheet1.
#
# Header cells follow the pattern:
#   "<Tagger> tag for <Judged> Scenario No. <N>"
#
# Tagger renames (applies to every header cell, since every header's
# tagger is one of these three):
#   Gemini-2.5-Pro      -> Grok-4-Fast
#   GPT-4o              -> GPT-5-mini
#   Claude-3.7-Sonnet   -> Mistral-Small-24b-2501
#
# Judged-model renames (applies only when the judged model name matches
# one of these four — other judged model names are left unchanged):
#   GPT-4o              -> GPT-5-mini
#   ChatGPT-4o          -> ChatGPT-5-mini
#   Gemini-2.5-Pro      -> Grok-4-Fast
#   Claude-3.7-Sonnet   -> Mistral-Small-24b-2501

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$taggerMap = @{
    "Gemini-2.5-Pro"    = "Grok-4-Fast";
    "GPT-4o"            = "GPT-5-mini";
    "Claude-3.7-Sonnet" = "Mistral-Small-24b-2501";
}

$judgedMap = @{
    "GPT-4o"            = "GPT-5-mini";
    "ChatGPT-4o"        = "ChatGPT-5-mini";
    "Gemini-2.5-Pro"    = "Grok-4-Fast";
    "Claude-3.7-Sonnet" = "Mistral-Small-24b-2501";
}

$lastCol = 757
$sep1 = " tag for "
$sep2 = " Scenario No. "

for ($col = 2; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $text = $cell.Value()

    $idx1 = $text.IndexOf($sep1)
    $tagger = $text.Substring(0, $idx1)
    $rest = $text.Substring($idx1 + $sep1.Length)

    $idx2 = $rest.IndexOf($sep2)
    $judged = $rest.Substring(0, $idx2)
    $num = $rest.Substring($idx2 + $sep2.Length)

    $newTagger = $tagger
    if ($taggerMap.ContainsKey($tagger)) {
        $newTagger = $taggerMap[$tagger]
    }

    $newJudged = $judged
    if ($judgedMap.ContainsKey($judged)) {
        $newJudged = $judgedMap[$judged]
    }

    $newText = $newTagger + $sep1 + $newJudged + $sep2 + $num
    $cell.Value = $newText
}
